$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they remain text like the original cells.
$textCells = @("D5","D6","D8","D11","D12","D13","D14","D19","D20","D21","D22","D25","D26","D27","D30","D33","D35","D36","D39","D40","D41","D45","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated price (D) and volume/change (E) values per row
$ws.Range("D2").Value = "62.860.93"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.459.93"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "570.27"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "146.92"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "2.459.91"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "26.89"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D17").Value = "62.937.63"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "2.442.33"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  +6.60%  "
$ws.Range("D21").Value = "323.57"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  +12.62%  "
$ws.Range("D25").Value = "66.28"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").Value = "621.45"
$ws.Range("E26").Value = "  +11.02%  "
$ws.Range("D27").Value = "8.62"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +8.95%  "
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +5.80%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").Value = "0.142"
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").Value = "5.07"
$ws.Range("E35").Value = "  +6.35%  "
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "18.68"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "144.65"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("E43").Value = "  +14.85%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "147.78"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("D47").Value = "20.76"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("D48").Value = "0.0537"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "0.602"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("D51").Value = "0.0919"
$ws.Range("E51").Value = "  -0.46%  "
